{"js": "// Collapse the two split \"<id>...</id>\" runs (produced when the inner\n// identifier text had different run formatting than the surrounding\n// \"<id>\"/\"</id>\" tag text) back into a single run with the tag text's\n// formatting, while also renaming the ids:\n//   p005r_a2 -> p005r_2\n//   p005v_a1 -> p005v_1\n\nconst body = context.document.body;\n\nconst replacements = [\n  [\"<id>p005r_a2</id>\", \"<id>p005r_2</id>\"],\n  [\"<id>p005v_a1</id>\", \"<id>p005v_1</id>\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Collapse the two split \"<id>...</id>\" runs (the inner identifier text had\n# different run formatting than the surrounding \"<id>\"/\"</id>\" tag text) back\n# into a single run with the tag text's formatting, while also renaming the\n# ids:\n#   p005r_a2 -> p005r_2\n#   p005v_a1 -> p005v_1\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\"<id>p005r_a2</id>\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"<id>p005r_2</id>\", $wdReplaceOne) | Out-Null\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"<id>p005v_a1</id>\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"<id>p005v_1</id>\", $wdReplaceOne) | Out-Null\n"}
